$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 2 (hunk @@ -727,25 +727,19 @@)
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()

# row 64 (hunk @@ -3816,25 +3810,25 @@)
$ws.Range("H64").Value = 4734.5
$ws.Range("I64").Value = 4566.6665
$ws.Range("J64").Value = 4860.375
$ws.Range("K64").Value = 4566.6665
$ws.Range("L64").Value = 4860.375
$ws.Range("M64").Value = -4318.6665
$ws.Range("N64").Value = -5356.375

# row 67 (hunk @@ -3966,25 +3960,25 @@)
$ws.Range("H67").Value = 4734.5
$ws.Range("I67").Value = 4566.6665
$ws.Range("J67").Value = 4860.375
$ws.Range("K67").Value = 4566.6665
$ws.Range("L67").Value = 4860.375
$ws.Range("M67").Value = -3708.6665
$ws.Range("N67").Value = -6576.375

# row 74 (hunk @@ -4318,25 +4312,25 @@)
$ws.Range("H74").Value = 5470.353
$ws.Range("I74").Value = 5399.7334
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 5399.7334
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = -4463.7334
$ws.Range("N74").Value = -7872

# row 77 (hunk @@ -4468,25 +4462,25 @@)
$ws.Range("H77").Value = 5470.353
$ws.Range("I77").Value = 5399.7334
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 26998.667
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = -22318.667
$ws.Range("N77").Value = -39360

# row 80 (hunk @@ -4618,25 +4612,25 @@)
$ws.Range("H80").Value = 955.6177
$ws.Range("I80").Value = 500.3125
$ws.Range("J80").Value = 1360.3334
$ws.Range("K80").Value = 1500.9375
$ws.Range("L80").Value = 4081.0002
$ws.Range("M80").Value = -502.9375
$ws.Range("N80").Value = -6077.0002

# row 83 (hunk @@ -4771,25 +4765,25 @@)
$ws.Range("H83").Value = 955.6177
$ws.Range("I83").Value = 500.3125
$ws.Range("J83").Value = 1360.3334
$ws.Range("K83").Value = 4502.8125
$ws.Range("L83").Value = 12243.0006
$ws.Range("M83").Value = 489.1875
$ws.Range("N83").Value = -22227.0006

$ws = $wb.Worksheets.Item("ARM")
# row 61 (hunk @@ -10719,25 +10713,25 @@)
$ws.Range("H61").Value = 5339.4326
$ws.Range("I61").Value = 5838.6333
$ws.Range("J61").Value = 3200
$ws.Range("K61").Value = 5838.6333
$ws.Range("L61").Value = 3200
$ws.Range("M61").Value = -5626.6333
$ws.Range("N61").Value = -3624

# row 74 (hunk @@ -11359,25 +11353,25 @@)
$ws.Range("H74").Value = 3083.102
$ws.Range("I74").Value = 780.7059
$ws.Range("J74").Value = 8301.866
$ws.Range("K74").Value = 780.7059
$ws.Range("L74").Value = 8301.866
$ws.Range("M74").Value = 93.29409999999996
$ws.Range("N74").Value = -10049.866

# row 77 (hunk @@ -11503,25 +11497,25 @@)
$ws.Range("H77").Value = 3083.102
$ws.Range("I77").Value = 780.7059
$ws.Range("J77").Value = 8301.866
$ws.Range("K77").Value = 3903.5295
$ws.Range("L77").Value = 41509.33
$ws.Range("M77").Value = 464.4704999999999
$ws.Range("N77").Value = -50245.33

# row 88 (hunk @@ -12039,25 +12033,25 @@)
$ws.Range("H88").Value = 14312.272
$ws.Range("I88").Value = 2000
$ws.Range("J88").Value = 21347.857
$ws.Range("K88").Value = 2000
$ws.Range("L88").Value = 21347.857
$ws.Range("M88").Value = -1594
$ws.Range("N88").Value = -22159.857

# row 91 (hunk @@ -12183,25 +12177,25 @@)
$ws.Range("H91").Value = 14312.272
$ws.Range("I91").Value = 2000
$ws.Range("J91").Value = 21347.857
$ws.Range("K91").Value = 2000
$ws.Range("L91").Value = 21347.857
$ws.Range("M91").Value = -596
$ws.Range("N91").Value = -24155.857

# row 122 (hunk @@ -13702,25 +13696,25 @@)
$ws.Range("H122").Value = 46155110
$ws.Range("I122").Value = 63159030
$ws.Range("J122").Value = 1619.1428
$ws.Range("K122").Value = 189477090
$ws.Range("L122").Value = 4857.428400000001
$ws.Range("M122").Value = -189474640
$ws.Range("N122").Value = -9757.428400000001

# row 136 (hunk @@ -14391,25 +14385,25 @@)
$ws.Range("H136").Value = 5339.4326
$ws.Range("I136").Value = 5838.6333
$ws.Range("J136").Value = 3200
$ws.Range("K136").Value = 17515.8999
$ws.Range("L136").Value = 9600
$ws.Range("M136").Value = -14965.8999
$ws.Range("N136").Value = -14700

$ws = $wb.Worksheets.Item("BSM")
# row 86 (hunk @@ -18856,25 +18850,25 @@)
$ws.Range("H86").Value = 1270.1
$ws.Range("I86").Value = 1250.1666
$ws.Range("J86").Value = 1300
$ws.Range("K86").Value = 1250.1666
$ws.Range("L86").Value = 1300
$ws.Range("M86").Value = -127.1666
$ws.Range("N86").Value = -3546

# row 89 (hunk @@ -19003,25 +18997,25 @@)
$ws.Range("H89").Value = 1270.1
$ws.Range("I89").Value = 1250.1666
$ws.Range("J89").Value = 1300
$ws.Range("K89").Value = 6250.833000000001
$ws.Range("L89").Value = 6500
$ws.Range("M89").Value = -634.8330000000005
$ws.Range("N89").Value = -17732

# row 99 (hunk @@ -19493,25 +19487,25 @@)
$ws.Range("H99").Value = 1074.9286
$ws.Range("I99").Value = 631.7273
$ws.Range("J99").Value = 2700
$ws.Range("K99").Value = 631.7273
$ws.Range("L99").Value = 2700
$ws.Range("M99").Value = 866.2727
$ws.Range("N99").Value = -5696

# row 134 (hunk @@ -21193,25 +21187,25 @@)
$ws.Range("H134").Value = 20604588
$ws.Range("I134").Value = 22244786
$ws.Range("J134").Value = 16668118
$ws.Range("K134").Value = 66734358
$ws.Range("L134").Value = 50004354
$ws.Range("M134").Value = -66731823
$ws.Range("N134").Value = -50009424

$ws = $wb.Worksheets.Item("CRP")
# row 31 (hunk @@ -23094,25 +23088,25 @@)
$ws.Range("H31").Value = 6648.4424
$ws.Range("I31").Value = 1233.72
$ws.Range("J31").Value = 10408.667
$ws.Range("K31").Value = 1233.72
$ws.Range("L31").Value = 10408.667
$ws.Range("M31").Value = -938.72
$ws.Range("N31").Value = -10998.667

# row 34 (hunk @@ -23247,25 +23241,25 @@)
$ws.Range("H34").Value = 6648.4424
$ws.Range("I34").Value = 1233.72
$ws.Range("J34").Value = 10408.667
$ws.Range("K34").Value = 1233.72
$ws.Range("L34").Value = 10408.667
$ws.Range("M34").Value = -1031.72
$ws.Range("N34").Value = -10812.667

# row 58 (hunk @@ -24405,25 +24399,25 @@)
$ws.Range("H58").Value = 2862595.5
$ws.Range("I58").Value = 5495325.5
$ws.Range("J58").Value = 10471.125
$ws.Range("K58").Value = 5495325.5
$ws.Range("L58").Value = 10471.125
$ws.Range("M58").Value = -5495122.5
$ws.Range("N58").Value = -10877.125

# row 133 (hunk @@ -28053,22 +28047,22 @@)
$ws.Range("H133").Value = 48000
$ws.Range("J133").Value = 48000
$ws.Range("L133").Value = 48000
$ws.Range("N133").Value = -53060

# row 136 (hunk @@ -28203,25 +28197,25 @@)
$ws.Range("H136").Value = 2862595.5
$ws.Range("I136").Value = 5495325.5
$ws.Range("J136").Value = 10471.125
$ws.Range("K136").Value = 16485976.5
$ws.Range("L136").Value = 31413.375
$ws.Range("M136").Value = -16483426.5
$ws.Range("N136").Value = -36513.375

$ws = $wb.Worksheets.Item("GSM")
# row 22 (hunk @@ -36825,19 +36819,25 @@)
$ws.Range("H22").Value = 2750
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = 29
$ws.Range("N22").Value = -6058

# row 122 (hunk @@ -41692,25 +41692,25 @@)
$ws.Range("H122").Value = 250005730
$ws.Range("I122").Value = 500000740
$ws.Range("J122").Value = 10700
$ws.Range("K122").Value = 1500002220
$ws.Range("L122").Value = 32100
$ws.Range("M122").Value = -1499999770
$ws.Range("N122").Value = -37000

$ws = $wb.Worksheets.Item("LTW")
# row 46 (hunk @@ -44940,22 +44940,22 @@)
$ws.Range("H46").Value = 21420
$ws.Range("I46").Value = 2550
$ws.Range("K46").Value = 2550
$ws.Range("M46").Value = -2362

# row 100 (hunk @@ -47565,25 +47565,25 @@)
$ws.Range("H100").Value = 2536.5
$ws.Range("I100").Value = 1671.5
$ws.Range("J100").Value = 3401.5
$ws.Range("K100").Value = 1671.5
$ws.Range("L100").Value = 3401.5
$ws.Range("M100").Value = -1130.5
$ws.Range("N100").Value = -4483.5

# row 122 (hunk @@ -48631,22 +48631,22 @@)
$ws.Range("H122").Value = 35717500
$ws.Range("J122").Value = 35717500
$ws.Range("L122").Value = 107152500
$ws.Range("N122").Value = -107157400

# row 132 (hunk @@ -49121,25 +49121,25 @@)
$ws.Range("H132").Value = 6266107
$ws.Range("I132").Value = 27782404
$ws.Range("J132").Value = 2167764.8
$ws.Range("K132").Value = 83347212
$ws.Range("L132").Value = 6503294.399999999
$ws.Range("M132").Value = -83344682
$ws.Range("N132").Value = -6508354.399999999

# row 136 (hunk @@ -49320,25 +49320,25 @@)
$ws.Range("H136").Value = 732326.9399999999
$ws.Range("I136").Value = 4223.9355
$ws.Range("J136").Value = 1986282.1
$ws.Range("K136").Value = 12671.8065
$ws.Range("L136").Value = 5958846.300000001
$ws.Range("M136").Value = -10121.8065
$ws.Range("N136").Value = -5963946.300000001

$ws = $wb.Worksheets.Item("WVR")
# row 136 (hunk @@ -56226,25 +56226,25 @@)
$ws.Range("H136").Value = 13115109
$ws.Range("I136").Value = 8051963.5
$ws.Range("J136").Value = 20000988
$ws.Range("K136").Value = 24155890.5
$ws.Range("L136").Value = 60002964
$ws.Range("M136").Value = -24153340.5
$ws.Range("N136").Value = -60008064
